$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1: 99.65 -> 0M
$tbl.Cell(1, 1).Range.Text = "0M"
# Row 2: 0.07 -> 0M
$tbl.Cell(2, 1).Range.Text = "0M"
# Row 3: 19 -> 0M
$tbl.Cell(3, 1).Range.Text = "0M"
# Row 4: 110 -> 310
$tbl.Cell(4, 1).Range.Text = "310"
# Row 5: 0.00006 -> 0.00001
$tbl.Cell(5, 1).Range.Text = "0.00001"
# Row 6: 0.00045 -> 0.00048
$tbl.Cell(6, 1).Range.Text = "0.00048"
# Row 7: 0.00015 -> 0.00017
$tbl.Cell(7, 1).Range.Text = "0.00017"
# Row 8 (0.00004) unchanged
# Row 9: 0.00019 -> 0.00032
$tbl.Cell(9, 1).Range.Text = "0.00032"
# Row 10: 0.00022 -> 0.00037
$tbl.Cell(10, 1).Range.Text = "0.00037"
# Row 11: 0.00027 -> 0.00040
$tbl.Cell(11, 1).Range.Text = "0.00040"
# Row 12: 0.02357 -> 0.06581
$tbl.Cell(12, 1).Range.Text = "0.06581"

# Row 44: collapse multi-value tabbed run down to single value "99.65"
$tbl.Cell(44, 1).Range.Text = "99.65"
# Row 45: collapse multi-value tabbed run down to single value "0.07"
$tbl.Cell(45, 1).Range.Text = "0.07"
# Row 46: collapse multi-value tabbed run down to single value "19"
$tbl.Cell(46, 1).Range.Text = "19"
